# Update row 8 (ano/ano_obj = 2025) metrics in the recurrence metrics sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 941
$ws.Range("D8").Value = 156
$ws.Range("E8").Value = 785
$ws.Range("F8").Value = 6.398687448728466
$ws.Range("G8").Value = 83.42189160467588
$ws.Range("H8").Value = 16.57810839532412
